# "finishing first line graph"
# Updates one existing data point and appends a new "TFrench Cities" data
# series (rows 34-41) to the Tuberculosis Towns worksheet, then restores
# the user's on-screen selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct an existing Rural France value (row 28, 1908) ---
$ws.Range("B28").Value = 1.62

# --- Append the new "TFrench Cities" series (rows 34-41) ---
$newSeries = @(
    @(1906, 3.25),
    @(1907, 3.32),
    @(1908, 3.25),
    @(1909, 3.13),
    @(1910, 3.14),
    @(1911, 3.15),
    @(1912, 3.03),
    @(1913, 3.07)
)

$row = 34
foreach ($point in $newSeries) {
    $ws.Cells.Item($row, 1).Value = $point[0]
    $ws.Cells.Item($row, 2).Value = $point[1]
    $ws.Cells.Item($row, 3).Value = "TFrench Cities"
    $row = $row + 1
}

# --- Restore the selection left by the author ---
$ws.Range("E27").Select() | Out-Null
